$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 588.3333
$ws.Range("I4").Value = 652
$ws.Range("J4").Value = 270
$ws.Range("K4").Value = 652
$ws.Range("L4").Value = 270
$ws.Range("M4").Value = -538
$ws.Range("H17").Value = 1964.4333
$ws.Range("J17").Value = 1964.4333
$ws.Range("L17").Value = 5893.2999
$ws.Range("N17").Value = -6229.2999
$ws.Range("H18").Value = 1790
$ws.Range("I18").Value = 1790
$ws.Range("K18").Value = 1790
$ws.Range("M18").Value = -1506
$ws.Range("H112").Value = 58825460
$ws.Range("J112").Value = 2051.25
$ws.Range("L112").Value = 6153.75
$ws.Range("N112").Value = -8369.75
$ws.Range("H129").Value = 891.8679
$ws.Range("J129").Value = 945.2766
$ws.Range("L129").Value = 2835.8298
$ws.Range("N129").Value = -12835.8298
$ws.Range("H137").Value = 2224128.8
$ws.Range("I137").Value = 3572983.8
$ws.Range("J137").Value = 2485.353
$ws.Range("K137").Value = 10718951.4
$ws.Range("L137").Value = 7456.059
$ws.Range("M137").Value = -10716401.4
$ws.Range("N137").Value = -12556.059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 57.4
$ws.Range("I5").Value = 57.4
$ws.Range("K5").Value = 57.4
$ws.Range("M5").Value = 54.6
$ws.Range("H32").Value = 3448846.5
$ws.Range("I32").Value = 3801110.2
$ws.Range("J32").Value = 4489.6665
$ws.Range("K32").Value = 3801110.2
$ws.Range("L32").Value = 4489.6665
$ws.Range("M32").Value = -3800823.2
$ws.Range("N32").Value = -5063.6665
$ws.Range("H45").Value = 2190.9092
$ws.Range("I45").Value = 3400
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 3400
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -3023
$ws.Range("N45").Value = -2254
$ws.Range("H61").Value = 100201270
$ws.Range("I61").Value = 125126090
$ws.Range("J61").Value = 502000
$ws.Range("K61").Value = 125126090
$ws.Range("L61").Value = 502000
$ws.Range("M61").Value = -125125878
$ws.Range("N61").Value = -502424
$ws.Range("H122").Value = 2742.4
$ws.Range("I122").Value = 2533.6
$ws.Range("K122").Value = 7600.799999999999
$ws.Range("M122").Value = -5150.799999999999
$ws.Range("H132").Value = 41022.117
$ws.Range("I132").Value = 28530.486
$ws.Range("J132").Value = 74035.71000000001
$ws.Range("K132").Value = 85591.458
$ws.Range("L132").Value = 222107.13
$ws.Range("M132").Value = -83061.458
$ws.Range("N132").Value = -227167.13
$ws.Range("H136").Value = 100201270
$ws.Range("I136").Value = 125126090
$ws.Range("J136").Value = 502000
$ws.Range("K136").Value = 375378270
$ws.Range("L136").Value = 1506000
$ws.Range("M136").Value = -375375720
$ws.Range("N136").Value = -1511100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 57.4
$ws.Range("I4").Value = 57.4
$ws.Range("K4").Value = 57.4
$ws.Range("M4").Value = 57.6
$ws.Range("H135").Value = 60500.4
$ws.Range("J135").Value = 60500.4
$ws.Range("L135").Value = 60500.4
$ws.Range("N135").Value = -70640.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 93.63636
$ws.Range("I7").Value = 68.75
$ws.Range("K7").Value = 68.75
$ws.Range("M7").Value = 44.25
$ws.Range("H122").Value = 2028.2667
$ws.Range("I122").Value = 1468.6666
$ws.Range("J122").Value = 4266.6665
$ws.Range("K122").Value = 4405.9998
$ws.Range("L122").Value = 12799.9995
$ws.Range("M122").Value = -1955.9998
$ws.Range("N122").Value = -17699.9995
$ws.Range("H132").Value = 44814.957
$ws.Range("I132").Value = 1337.2
$ws.Range("J132").Value = 334666.66
$ws.Range("K132").Value = 4011.6
$ws.Range("L132").Value = 1003999.98
$ws.Range("M132").Value = -1481.6
$ws.Range("N132").Value = -1009059.98
$ws.Range("H133").Value = 41669.43
$ws.Range("J133").Value = 41669.43
$ws.Range("L133").Value = 41669.43
$ws.Range("N133").Value = -46729.43

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2400.4
$ws.Range("I80").Value = 1001
$ws.Range("K80").Value = 3003
$ws.Range("M80").Value = -11871.9999
$ws.Range("H83").Value = 2400.4
$ws.Range("I83").Value = 1001
$ws.Range("K83").Value = 9009
$ws.Range("M83").Value = -39359.9997
$ws.Range("H97").Value = 2011.7693
$ws.Range("I97").Value = 2914.7144
$ws.Range("J97").Value = 958.3333
$ws.Range("K97").Value = 8744.143199999999
$ws.Range("L97").Value = 2874.9999
$ws.Range("M97").Value = -8248.143199999999
$ws.Range("N97").Value = -3866.9999
$ws.Range("H107").Value = 923.386
$ws.Range("I107").Value = 482.14285
$ws.Range("J107").Value = 1180.7778
$ws.Range("K107").Value = 1446.42855
$ws.Range("L107").Value = 3542.3334
$ws.Range("M107").Value = 473.5714499999999
$ws.Range("N107").Value = -7382.3334
$ws.Range("H131").Value = 974.37933
$ws.Range("J131").Value = 1022.7692
$ws.Range("L131").Value = 3068.3076
$ws.Range("N131").Value = -13148.3076

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2119.75
$ws.Range("I102").Value = 2115.375
$ws.Range("K102").Value = 2115.375
$ws.Range("M102").Value = -493.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4597.8
$ws.Range("I7").Value = 3334.8
$ws.Range("J7").Value = 5860.8
$ws.Range("K7").Value = 3334.8
$ws.Range("L7").Value = 5860.8
$ws.Range("M7").Value = -3222.8
$ws.Range("N7").Value = -6084.8
$ws.Range("H22").Value = 26316354
$ws.Range("I22").Value = 450
$ws.Range("K22").Value = 450
$ws.Range("M22").Value = -155
$ws.Range("H27").Value = 26316354
$ws.Range("I27").Value = 450
$ws.Range("K27").Value = 450
$ws.Range("M27").Value = -343
$ws.Range("H74").Value = 3349399
$ws.Range("I74").Value = 10000197
$ws.Range("K74").Value = 10000197
$ws.Range("M74").Value = -9999199
$ws.Range("H77").Value = 3349399
$ws.Range("I77").Value = 10000197
$ws.Range("K77").Value = 30000591
$ws.Range("M77").Value = -29995599
$ws.Range("H122").Value = 3465.7144
$ws.Range("I122").Value = 3465.7144
$ws.Range("K122").Value = 10397.1432
$ws.Range("M122").Value = -7947.143199999999
$ws.Range("H126").Value = 4597.8
$ws.Range("I126").Value = 3334.8
$ws.Range("J126").Value = 5860.8
$ws.Range("K126").Value = 10004.4
$ws.Range("L126").Value = 17582.4
$ws.Range("M126").Value = -7534.400000000001
$ws.Range("N126").Value = -22522.4
$ws.Range("H136").Value = 76551.72
$ws.Range("I136").Value = 40369.883
$ws.Range("J136").Value = 148915.39
$ws.Range("K136").Value = 121109.649
$ws.Range("L136").Value = 446746.17
$ws.Range("M136").Value = -118559.649
$ws.Range("N136").Value = -451846.17
